$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 330, shifting existing rows 330:367 down to 331:368
$ws.Rows.Item(330).Insert()

# Populate the newly inserted row with the new weekly data point
$ws.Range("A330").Value = 10
$ws.Range("B330").Value = "Vega Modelo de Temuco"
$ws.Range("C330").Value = "La Araucanía"
$ws.Range("D330").Value = 45212
$ws.Range("E330").Value = 9
$ws.Range("F330").Value = 100114007
$ws.Range("G330").Value = "Jengibre"
$ws.Range("H330").Value = "Sin especificar"
$ws.Range("I330").Value = "Primera"
$ws.Range("J330").Value = 50
$ws.Range("K330").Value = 23000
$ws.Range("L330").Value = 23000
$ws.Range("M330").Value = 23000
$ws.Range("N330").Value = "$/caja 13 kilos"
$ws.Range("O330").Value = "Perú"
$ws.Range("P330").Value = 1769
$ws.Range("Q330").Value = 13
$ws.Range("R330").Value = "Hortaliza"
